$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''309.52'
$ws.Range('E2').Value = '''1.14%'
$ws.Range('D3').Value = '''35.57'
$ws.Range('E3').Value = '''-1.70%'
$ws.Range('D4').Value = '''5.124'
$ws.Range('E4').Value = '''1.46%'
$ws.Range('D5').Value = '''0.08186'
$ws.Range('E5').Value = '''3.59%'
$ws.Range('D6').Value = '''2.048'
$ws.Range('D7').Value = '''7.957'
$ws.Range('E7').Value = '''-0.50%'
$ws.Range('B8').Value = '''BTSEToken'
$ws.Range('C8').Value = '''https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D8').Value = '''2.896'
$ws.Range('E8').Value = '''8.85%'
$ws.Range('B9').Value = '''MXToken'
$ws.Range('C9').Value = '''https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D9').Value = '''0.9280'
$ws.Range('E9').Value = '''0.10%'
$ws.Range('B10').Value = '''LiechtensteinCryptoassetsExchange'
$ws.Range('C10').Value = '''https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D10').Value = '''0.1076'
$ws.Range('E10').Value = '''9.22%'
$ws.Range('B11').Value = '''WazirX'
$ws.Range('C11').Value = '''https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').Value = '''0.1918'
$ws.Range('E11').Value = '''2.71%'
$ws.Range('B12').Value = '''MandalaExchangeToken'
$ws.Range('C12').Value = '''https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').Value = '''0.09566'
$ws.Range('E12').Value = '''6.76%'
$ws.Range('B13').Value = '''BitrueCoin'
$ws.Range('C13').Value = '''https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').Value = '''0.03621'
$ws.Range('E13').Value = '''-3.54%'
$ws.Range('B14').Value = '''BitMartToken'
$ws.Range('C14').Value = '''https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').Value = '''0.09906'
$ws.Range('E14').Value = '''-0.11%'
$ws.Range('B15').Value = '''BitForexToken'
$ws.Range('C15').Value = '''https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').Value = '''0.001429'
$ws.Range('E15').Value = '''-0.93%'
$ws.Range('B16').Value = '''TigerCash'
$ws.Range('C16').Value = '''https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').Value = '''0.005696'
$ws.Range('E16').Value = '''-0.40%'
$ws.Range('B17').Value = '''LEO'
$ws.Range('C17').Value = '''https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').Value = '''3.470'
$ws.Range('E17').Value = '''0.14%'
$ws.Range('B18').Value = '''GateToken'
$ws.Range('C18').Value = '''https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D18').Value = '''4.134'
$ws.Range('E18').Value = '''-0.41%'
$ws.Range('D19').Value = '''0.3415'
$ws.Range('E19').Value = '''1.42%'
$ws.Range('D20').Value = '''0.1312'
$ws.Range('E20').Value = '''-0.55%'
$ws.Range('D21').Value = '''5.106'
$ws.Range('E21').Value = '''0.31%'
$ws.Range('D22').Value = '''0.2196'
$ws.Range('E22').Value = '''-2.42%'
$ws.Range('D23').Value = '''0.04544'
$ws.Range('E23').Value = '''-0.71%'
$ws.Range('E24').Value = '''-0.49%'
$ws.Range('D25').Value = '''0.004778'
$ws.Range('E25').Value = '''0.02%'
$ws.Range('D26').Value = '''0.0001253'
$ws.Range('E26').Value = '''-3.68%'
$ws.Range('D27').Value = '''0.0004459'
$ws.Range('E27').Value = '''-5.90%'
$ws.Range('D39').Value = '''0.01962'
$ws.Range('E39').Value = '''2.00%'
$ws.Range('D40').Value = '''0.04878'
$ws.Range('E40').Value = '''-1.02%'
$ws.Range('D41').Value = '''0.007687'
$ws.Range('E41').Value = '''-1.47%'
$ws.Range('D42').Value = '''0.009843'
$ws.Range('E42').Value = '''26.01%'
$ws.Range('D43').Value = '''0.1380'
$ws.Range('E43').Value = '''-0.79%'
$ws.Range('D44').Value = '''0.002120'
$ws.Range('E44').Value = '''-2.84%'
$ws.Range('D45').Value = '''0.01155'
$ws.Range('E45').Value = '''1.13%'
$ws.Range('D46').Value = '''0.00006516'
$ws.Range('E46').Value = '''6.00%'
$ws.Range('D47').Value = '''0.00000000752'
$ws.Range('E47').Value = '''0.16%'
$ws.Range('D48').Value = '''64.36'
$ws.Range('E48').Value = '''24.32%'
$ws.Range('E49').Value = '''-16.80%'
$ws.Range('D50').Value = '''0.00002105'
$ws.Range('E50').Value = '''0.16%'
$ws.Range('D51').Value = '''0.0002005'
$ws.Range('E51').Value = '''0.16%'

Write-Output "Applied all cell updates."
